$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Picture") to make room for "Barcode".
# This shifts the old B,C,D,E (Picture, Number_in_stock, Price, Description)
# right by one, to C,D,E,F - matching the diff.
$ws.Columns("B").Insert()

# Rename the old "Item_id" header (now still in A1) to "Name".
$ws.Range("A1").Value = "Name"

# New column B header.
$ws.Range("B1").Value = "Barcode"

# The Number_in_stock column (now column D after the insert) gets overwritten
# to a flat value of 1 for every data row.
$ws.Range("D2:D10").Value = 1

# Update the active selection to match the author's saved cursor position.
$ws.Range("G13").Select()

$wb.Save()
